# Nettlinx Ltd - Quarterly sheet: add "Exceptional items" column
#
# The "Quarterly" worksheet gets a brand-new column inserted right before
# the old column L ("P/l before tax"). The new column carries the label
# "Exceptional items" (row 1, plain-case header) / "Exceptional Items"
# (row 2, title-case header) and holds the exceptional-items figure that
# used to be silently folded into the old "P/l before tax" column. Only
# quarter Dec'20 (data row 28) actually had a non-zero exceptional item
# (0.41 = 0.73 - 0.32); every other quarter's value is blank/zero there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new blank column before the existing column L; everything
# from L onward (through the old T) shifts one column to the right
# (new M .. U), automatically carrying over values, formulas and styles.
$ws.Columns("L").Insert()

# Populate the header cells for the freshly-inserted column.
$ws.Range("L1").Value2 = "Exceptional items"
$ws.Range("L2").Value2 = "Exceptional Items"

# The only data row with a genuine exceptional-items figure: Dec '20.
$ws.Range("L28").Value2 = 0.41
